$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 3 roster change: Daniel Burcham -> Shelia Lowe
$ws.Range("A17").Value = "Shelia Lowe"

# Week-2 / Week-3 attendance/result updates (were placeholder "A")
$ws.Range("C3").Value = "DNP"
$ws.Range("D3").Value = "L"

$ws.Range("C4").Value = "W"
$ws.Range("D4").Value = "L"

$ws.Range("C5").Value = "L"
$ws.Range("D5").Value = "NA"

$ws.Range("C6").Value = "NA"
$ws.Range("D6").Value = "W"

$ws.Range("C7").Value = "L"

$ws.Range("C8").Value = "NA"
$ws.Range("D8").Value = "L"

$ws.Range("C9").Value = "L"
$ws.Range("D9").Value = "L"

$ws.Range("C10").Value = "L"
$ws.Range("D10").Value = "DNP"

$ws.Range("C15").Value = "L"
$ws.Range("D15").Value = "W"

$ws.Range("C16").Value = "NA"
$ws.Range("D16").Value = "L"

$ws.Range("B17").Value = "NA"
$ws.Range("C17").Value = "NA"
$ws.Range("D17").Value = "L"

$ws.Range("C18").Value = "W"

$ws.Range("C19").Value = "L"
$ws.Range("D19").Value = "DNP"

$ws.Range("C20").Value = "NA"
$ws.Range("D20").Value = "NA"

$ws.Range("C21").Value = "DNP"
$ws.Range("D21").Value = "L"

$ws.Range("C22").Value = "W"

# Selection as left by the author after editing
$ws.Range("E6").Select()
